# Powerpoint writer: consolidate text run nodes.
# Merge adjacent runs that were split purely because of a trailing/leading
# space (e.g. "Jesse" + " " + "Rosenthal" -> "Jesse " + "Rosenthal") by
# rewriting the leading substring in place; this preserves the other runs,
# their <a:rPr/> and the paragraph's <a:pPr/>.

$p = $ppt.ActivePresentation

# --- Slide 1, Subtitle "Jesse Rosenthal" ---------------------------------
# Text (after the two manual line breaks) is "Jesse Rosenthal":
#   pos 3-7 "Jesse", pos 8 " ", pos 9-17 "Rosenthal"
# Merge "Jesse" + " " into a single run "Jesse ".
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Characters(3, 6).Text = "Jesse "

# --- Slide 1 notes, "Some speaker notes" ---------------------------------
# Merge "Some" + " " into "Some ", and "speaker" + " " into "speaker ",
# leaving the trailing "notes" run untouched.
$notesRange = $slide1.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesRange.Text = "Some speaker notes"

# --- Slide 2, Title "A header" -------------------------------------------
# Merge "A" + " " into "A ".
$slide2 = $p.Slides.Item(2)
$title2 = $slide2.Shapes.Item(1)
$title2.TextFrame.TextRange.Characters(1, 2).Text = "A "
